# Auto-generated Excel COM-interop script.
#
# Commit message: "Updated symbol list on Mon Feb 13 16:16:34 UTC 2023 with
# GitHub Actions" -- a scheduled scraper refreshed the cryptocurrency table
# on Sheet1. This applies the resulting 136 cell-value changes (coin
# names/links shuffling down a slot, refreshed prices / 1h-volume deltas,
# and the "Hora" column bumping from 15 -> 16).
#
# Every touched cell in the source workbook is stored as literal TEXT
# (e.g. "287.48", "-9.34%", "16") rather than as a numeric/percentage
# value, even though the strings look numeric. Plain `Range.Value =
# "287.48"` would get auto-coerced by Excel into a real number (and a
# percentage string into a percent-formatted number), which would change
# both the stored type and the cell style. To avoid that, each cell is
# temporarily switched to the Text number format ("@") before the value
# is written, then its style is reset back to "Normal" so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @{Cell="D2"; Value="287.48"},
    @{Cell="E2"; Value="-9.34%"},
    @{Cell="G2"; Value="16"},
    @{Cell="D3"; Value="40.38"},
    @{Cell="E3"; Value="-1.86%"},
    @{Cell="G3"; Value="16"},
    @{Cell="D4"; Value="5.030"},
    @{Cell="E4"; Value="-3.38%"},
    @{Cell="G4"; Value="16"},
    @{Cell="D5"; Value="0.07318"},
    @{Cell="E5"; Value="-4.96%"},
    @{Cell="G5"; Value="16"},
    @{Cell="D6"; Value="4.282"},
    @{Cell="E6"; Value="-0.51%"},
    @{Cell="G6"; Value="16"},
    @{Cell="D7"; Value="1.544"},
    @{Cell="E7"; Value="-9.22%"},
    @{Cell="G7"; Value="16"},
    @{Cell="D8"; Value="0.9133"},
    @{Cell="E8"; Value="-3.64%"},
    @{Cell="G8"; Value="16"},
    @{Cell="D9"; Value="0.1202"},
    @{Cell="E9"; Value="-5.30%"},
    @{Cell="G9"; Value="16"},
    @{Cell="D10"; Value="0.1746"},
    @{Cell="E10"; Value="-4.65%"},
    @{Cell="G10"; Value="16"},
    @{Cell="D11"; Value="0.08671"},
    @{Cell="E11"; Value="-4.60%"},
    @{Cell="G11"; Value="16"},
    @{Cell="D12"; Value="0.04161"},
    @{Cell="E12"; Value="-1.71%"},
    @{Cell="G12"; Value="16"},
    @{Cell="E13"; Value="-0.13%"},
    @{Cell="G13"; Value="16"},
    @{Cell="D14"; Value="0.001276"},
    @{Cell="E14"; Value="-0.20%"},
    @{Cell="G14"; Value="16"},
    @{Cell="B15"; Value="CoinExToken"},
    @{Cell="C15"; Value="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"},
    @{Cell="D15"; Value="0.03848"},
    @{Cell="E15"; Value="-4.25%"},
    @{Cell="G15"; Value="16"},
    @{Cell="B16"; Value="TigerCash"},
    @{Cell="C16"; Value="https://coinranking.com/coin/6hIn06L2+tigercash-tch"},
    @{Cell="D16"; Value="0.005782"},
    @{Cell="E16"; Value="-1.65%"},
    @{Cell="G16"; Value="16"},
    @{Cell="B17"; Value="LEO"},
    @{Cell="C17"; Value="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"},
    @{Cell="D17"; Value="3.400"},
    @{Cell="E17"; Value="1.35%"},
    @{Cell="G17"; Value="16"},
    @{Cell="B18"; Value="BTSEToken"},
    @{Cell="C18"; Value="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"},
    @{Cell="D18"; Value="2.397"},
    @{Cell="E18"; Value="-1.16%"},
    @{Cell="G18"; Value="16"},
    @{Cell="B19"; Value="BitpandaEcosystemToken"},
    @{Cell="C19"; Value="https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"},
    @{Cell="D19"; Value="0.3273"},
    @{Cell="E19"; Value="-2.47%"},
    @{Cell="G19"; Value="16"},
    @{Cell="B20"; Value="MCDex"},
    @{Cell="C20"; Value="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"},
    @{Cell="D20"; Value="7.558"},
    @{Cell="E20"; Value="1.34%"},
    @{Cell="G20"; Value="16"},
    @{Cell="B21"; Value="ProBitToken"},
    @{Cell="C21"; Value="https://coinranking.com/coin/lQP4d6T2+probittoken-prob"},
    @{Cell="D21"; Value="0.1342"},
    @{Cell="E21"; Value="-0.70%"},
    @{Cell="G21"; Value="16"},
    @{Cell="B22"; Value="ZBToken"},
    @{Cell="C22"; Value="https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"},
    @{Cell="D22"; Value="0.2884"},
    @{Cell="E22"; Value="3.63%"},
    @{Cell="G22"; Value="16"},
    @{Cell="E23"; Value="0.31%"},
    @{Cell="G23"; Value="16"},
    @{Cell="D24"; Value="0.003785"},
    @{Cell="E24"; Value="-10.66%"},
    @{Cell="G24"; Value="16"},
    @{Cell="D25"; Value="0.0001282"},
    @{Cell="E25"; Value="0.83%"},
    @{Cell="G25"; Value="16"},
    @{Cell="D26"; Value="0.0003727"},
    @{Cell="E26"; Value="-95.04%"},
    @{Cell="G26"; Value="16"},
    @{Cell="G27"; Value="16"},
    @{Cell="G28"; Value="16"},
    @{Cell="G29"; Value="16"},
    @{Cell="G30"; Value="16"},
    @{Cell="G31"; Value="16"},
    @{Cell="G32"; Value="16"},
    @{Cell="G33"; Value="16"},
    @{Cell="G34"; Value="16"},
    @{Cell="G35"; Value="16"},
    @{Cell="G36"; Value="16"},
    @{Cell="G37"; Value="16"},
    @{Cell="D38"; Value="0.02330"},
    @{Cell="E38"; Value="-8.10%"},
    @{Cell="G38"; Value="16"},
    @{Cell="D39"; Value="0.05008"},
    @{Cell="E39"; Value="-6.04%"},
    @{Cell="G39"; Value="16"},
    @{Cell="D40"; Value="0.007684"},
    @{Cell="E40"; Value="-1.71%"},
    @{Cell="G40"; Value="16"},
    @{Cell="E41"; Value="163.47%"},
    @{Cell="G41"; Value="16"},
    @{Cell="D42"; Value="0.1272"},
    @{Cell="E42"; Value="-3.20%"},
    @{Cell="G42"; Value="16"},
    @{Cell="D43"; Value="0.007381"},
    @{Cell="E43"; Value="0.31%"},
    @{Cell="G43"; Value="16"},
    @{Cell="D44"; Value="0.007508"},
    @{Cell="E44"; Value="-0.69%"},
    @{Cell="G44"; Value="16"},
    @{Cell="D45"; Value="0.3113"},
    @{Cell="E45"; Value="-1.53%"},
    @{Cell="G45"; Value="16"},
    @{Cell="D46"; Value="0.00006486"},
    @{Cell="E46"; Value="-3.31%"},
    @{Cell="G46"; Value="16"},
    @{Cell="E47"; Value="0.15%"},
    @{Cell="G47"; Value="16"},
    @{Cell="E48"; Value="14.12%"},
    @{Cell="G48"; Value="16"},
    @{Cell="G49"; Value="16"},
    @{Cell="E50"; Value="0.15%"},
    @{Cell="G50"; Value="16"},
    @{Cell="D51"; Value="0.0002002"},
    @{Cell="E51"; Value="0.15%"},
    @{Cell="G51"; Value="16"}
)

foreach ($update in $cellUpdates) {
    $rng = $ws.Range($update.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $update.Value
    $rng.Style = "Normal"
}
